$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "73.302.29"
$ws.Range("E2").Value = "  +1.86%  "

$ws.Range("D3").Value = "4.053.77"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("E4").Value = "  -0.06%  "

Set-TextValue "D5" "575.10"
$ws.Range("E5").Value = "  +7.92%  "

Set-TextValue "D6" "152.61"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").Value = "4.049.55"
$ws.Range("E7").Value = "  +1.29%  "

Set-TextValue "D8" "0.697"
$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  +2.88%  "

$ws.Range("E11").Value = "  +0.59%  "

Set-TextValue "D12" "54.06"
$ws.Range("E12").Value = "  +13.98%  "

Set-TextValue "D13" "0.0000328"
$ws.Range("E13").Value = "  +0.92%  "

Set-TextValue "D14" "11.25"
$ws.Range("E14").Value = "  +5.71%  "

$ws.Range("D15").Value = "4.708.72"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "4.055.11"
$ws.Range("E16").Value = "  +1.14%  "

Set-TextValue "D17" "14.40"
$ws.Range("E17").Value = "  +3.39%  "

Set-TextValue "D18" "20.91"
$ws.Range("E18").Value = "  +2.12%  "

Set-TextValue "D19" "1.23"
$ws.Range("E19").Value = "  +3.38%  "

$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").Value = "73.250.43"
$ws.Range("E21").Value = "  +1.92%  "

Set-TextValue "D22" "445.95"
$ws.Range("E22").Value = "  +4.71%  "

$ws.Range("E23").Value = "  +8.50%  "

Set-TextValue "D24" "98.49"
$ws.Range("E24").Value = "  +0.74%  "

Set-TextValue "D25" "3.58"
$ws.Range("E25").Value = "  +2.40%  "

Set-TextValue "D26" "14.75"
$ws.Range("E26").Value = "  +3.05%  "

Set-TextValue "D27" "4.27"
$ws.Range("E27").Value = "  +18.54%  "

Set-TextValue "D28" "11.52"
$ws.Range("E28").Value = "  +3.15%  "

Set-TextValue "D29" "11.10"
$ws.Range("E29").Value = "  +4.04%  "

$ws.Range("E30").Value = "  +1.86%  "

Set-TextValue "D31" "37.27"
$ws.Range("E31").Value = "  +1.52%  "

Set-TextValue "D32" "7.90"
$ws.Range("E32").Value = "  +10.83%  "

$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D33" "13.76"
$ws.Range("E33").Value = "  +3.11%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.134"
$ws.Range("E34").Value = "  +4.38%  "

Set-TextValue "D35" "688.37"
$ws.Range("E35").Value = "  +2.62%  "

Set-TextValue "D36" "48.66"
$ws.Range("E36").Value = "  +14.17%  "

Set-TextValue "D37" "68.20"
$ws.Range("E37").Value = "  +3.74%  "

$ws.Range("D38").Value = "0.0₃0904"
$ws.Range("E38").Value = "  +9.75%  "

$ws.Range("E39").Value = "  +5.22%  "

$ws.Range("E40").Value = "  -1.80%  "

Set-TextValue "D41" "11.34"
$ws.Range("E41").Value = "  +17.39%  "

Set-TextValue "D42" "3.38"
$ws.Range("E42").Value = "  -1.37%  "

Set-TextValue "D43" "0.998"
$ws.Range("E43").Value = "  +0.03%  "

Set-TextValue "D44" "3.33"
$ws.Range("E44").Value = "  +1.48%  "

Set-TextValue "D45" "0.0496"
$ws.Range("E45").Value = "  +1.93%  "

Set-TextValue "D46" "1.00"
$ws.Range("E46").Value = "  +0.13%  "

Set-TextValue "D47" "0.152"
$ws.Range("E47").Value = "  +1.50%  "

Set-TextValue "D48" "2.77"
$ws.Range("E48").Value = "  +5.51%  "

Set-TextValue "D49" "2.25"
$ws.Range("E49").Value = "  +12.97%  "

Set-TextValue "D50" "3.54"
$ws.Range("E50").Value = "  +7.83%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D51" "3.07"
$ws.Range("E51").Value = "  +3.24%  "
